$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.143.41"
$ws.Range("E2").Value = "  -3.05%  "

$ws.Range("D3").Value = "1.711.67"
$ws.Range("E3").Value = "  -3.51%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.83"
$ws.Range("E5").Value = "  -5.95%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4764"
$ws.Range("E7").Value = "  +5.87%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3439"
$ws.Range("E8").Value = "  -3.31%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.13"
$ws.Range("E9").Value = "  +0.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07279"
$ws.Range("E10").Value = "  -2.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.040"
$ws.Range("E11").Value = "  -6.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.82"
$ws.Range("E13").Value = "  -5.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.851"
$ws.Range("E14").Value = "  -3.30%  "

$ws.Range("D15").Value = "1.709.77"
$ws.Range("E15").Value = "  -3.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.835"
$ws.Range("E16").Value = "  -5.86%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.62"
$ws.Range("E17").Value = "  -5.76%  "

$ws.Range("E18").Value = "  -2.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06354"
$ws.Range("E19").Value = "  -1.41%  "

$ws.Range("E20").Value = "  +0.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.47"
$ws.Range("E21").Value = "  -3.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.609"
$ws.Range("E22").Value = "  -3.10%  "

$ws.Range("D23").Value = "27.174.76"
$ws.Range("E23").Value = "  -2.98%  "

$ws.Range("E24").Value = "  -4.59%  "

$ws.Range("E25").Value = "  -1.95%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.86"
$ws.Range("E26").Value = "  -5.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.65"
$ws.Range("E27").Value = "  -3.67%  "

$ws.Range("D28").Value = "1.905.41"
$ws.Range("E28").Value = "  -3.54%  "

$ws.Range("E29").Value = "  -4.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.03"
$ws.Range("E30").Value = "  -3.72%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.013"
$ws.Range("E31").Value = "  -8.72%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09263"
$ws.Range("E32").Value = "  +0.50%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.590"
$ws.Range("E33").Value = "  -2.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.291"
$ws.Range("E34").Value = "  -7.35%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02193"
$ws.Range("E35").Value = "  -4.39%  "

$ws.Range("E36").Value = "  -5.08%  "

$ws.Range("E37").Value = "  -7.06%  "

$ws.Range("E38").Value = "  -4.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.751"
$ws.Range("E39").Value = "  -4.77%  "

$ws.Range("E40").Value = "  +0.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.000"
$ws.Range("E41").Value = "  +0.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5918"
$ws.Range("E42").Value = "  -6.34%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.114"
$ws.Range("E43").Value = "  -6.26%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.493"
$ws.Range("E44").Value = "  -5.34%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.63"
$ws.Range("E45").Value = "  -5.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.565"
$ws.Range("E46").Value = "  -5.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5615"
$ws.Range("E47").Value = "  -4.76%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "118.26"
$ws.Range("E48").Value = "  -3.61%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.837"
$ws.Range("E49").Value = "  -6.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06633"
$ws.Range("E50").Value = "  -3.75%  "

$ws.Range("E51").Value = "  -5.06%  "
